$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel PasteSpecial constants
$xlPasteValues = -4163
$xlPasteFormats = -4122

# The sheet currently only has data rows for 15-01-2024..22-01-2024 (rows 2-9).
# The target data runs from 01-01-2024..25-01-2024 (25 rows, rows 2-26), so
# rows 10-26 are brand new. Give them the same bold/centered/bordered style
# that column A already uses (copied from the existing A9 cell) before
# filling in any content.
$ws.Range("A9").Copy()
$ws.Range("A10:A26").PasteSpecial($xlPasteFormats)

# Write date text values via a scratch column (Z) so that Excel does not
# auto-convert the "dd-mm-yyyy" looking strings into real date serials when
# they are assigned to column A. Each Z cell gets a formula that evaluates
# to the literal text, then the computed text results are copied and
# pasted as values (not formulas) into column A, after which the scratch
# column is cleared again.
$ws.Range("Z2").Formula = "=""01-01-2024"""
$ws.Range("Z3").Formula = "=""02-01-2024"""
$ws.Range("Z4").Formula = "=""03-01-2024"""
$ws.Range("Z5").Formula = "=""04-01-2024"""
$ws.Range("Z6").Formula = "=""05-01-2024"""
$ws.Range("Z7").Formula = "=""06-01-2024"""
$ws.Range("Z8").Formula = "=""07-01-2024"""
$ws.Range("Z9").Formula = "=""08-01-2024"""
$ws.Range("Z10").Formula = "=""09-01-2024"""
$ws.Range("Z11").Formula = "=""10-01-2024"""
$ws.Range("Z12").Formula = "=""11-01-2024"""
$ws.Range("Z13").Formula = "=""12-01-2024"""
$ws.Range("Z14").Formula = "=""13-01-2024"""
$ws.Range("Z15").Formula = "=""14-01-2024"""
$ws.Range("Z16").Formula = "=""15-01-2024"""
$ws.Range("Z17").Formula = "=""16-01-2024"""
$ws.Range("Z18").Formula = "=""17-01-2024"""
$ws.Range("Z19").Formula = "=""18-01-2024"""
$ws.Range("Z20").Formula = "=""19-01-2024"""
$ws.Range("Z21").Formula = "=""20-01-2024"""
$ws.Range("Z22").Formula = "=""21-01-2024"""
$ws.Range("Z23").Formula = "=""22-01-2024"""
$ws.Range("Z24").Formula = "=""23-01-2024"""
$ws.Range("Z25").Formula = "=""24-01-2024"""
$ws.Range("Z26").Formula = "=""25-01-2024"""
$ws.Range("Z2:Z26").Copy()
$ws.Range("A2:A26").PasteSpecial($xlPasteValues)
$ws.Range("Z2:Z26").ClearContents()

# Write the numeric energy [kWh] values in column B
$ws.Range("B2").Value = 0.2640727228724609
$ws.Range("B3").Value = 0.0002728733781852322
$ws.Range("B4").Value = 0.01318142747795677
$ws.Range("B5").Value = 0.02769346615424148
$ws.Range("B6").Value = 0.1045413447121037
$ws.Range("B7").Value = 0.03146386253696185
$ws.Range("B8").Value = 0.003810482642393739
$ws.Range("B9").Value = 0.5208876517689147
$ws.Range("B10").Value = 0.9372133312270073
$ws.Range("B11").Value = 0.8436036806890808
$ws.Range("B12").Value = 0.8617712033086313
$ws.Range("B13").Value = 0.02167194237400416
$ws.Range("B14").Value = 0.02033969313737993
$ws.Range("B15").Value = 0.01212794262982757
$ws.Range("B16").Value = 0.01173135699545329
$ws.Range("B17").Value = 0.08389218667816668
$ws.Range("B18").Value = 0.008029626214607797
$ws.Range("B19").Value = 0.001148125312265952
$ws.Range("B20").Value = 0.3920495730652868
$ws.Range("B21").Value = 1.16301139171748
$ws.Range("B22").Value = 0.4922478884896943
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0.4999612913260334
$ws.Range("B25").Value = 0.2215556335889916
$ws.Range("B26").Value = 0.3556320059476869

Write-Host "A2=$($ws.Range('A2').Text) B2=$($ws.Range('B2').Text)"
Write-Host "A26=$($ws.Range('A26').Text) B26=$($ws.Range('B26').Text)"
